# Insert a new weekly-price record as row 56, pushing every existing
# record (old rows 56..127) down by one (new rows 57..128).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Insert()

$ws.Cells.Item(56, 1).Value2 = 8
$ws.Cells.Item(56, 2).Value2 = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(56, 3).Value2 = 'Coquimbo'
$ws.Cells.Item(56, 4).Value2 = 44848
$ws.Cells.Item(56, 5).Value2 = 4
$ws.Cells.Item(56, 6).Value2 = 100112052
$ws.Cells.Item(56, 7).Value2 = 'Albahaca'
$ws.Cells.Item(56, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(56, 9).Value2 = 'Primera'
$ws.Cells.Item(56, 10).Value2 = 1100
$ws.Cells.Item(56, 11).Value2 = 4000
$ws.Cells.Item(56, 12).Value2 = 4500
$ws.Cells.Item(56, 13).Value2 = 4250
$ws.Cells.Item(56, 14).Value2 = '$/paquete'
$ws.Cells.Item(56, 15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(56, 16).Value2 = 4250
$ws.Cells.Item(56, 17).Value2 = 1
$ws.Cells.Item(56, 18).Value2 = 'Hortaliza'
